# Auto-generated script applying crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates -- force text storage so values are not reinterpreted as numbers/dates,
# matching the original inlineStr text cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.961.44'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.04'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.02'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.256'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.62'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.870.33'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.28'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.644.58'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.88'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.918.46'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.12'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.93'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.44'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.904'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.132.07'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.540'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.35'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.796'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.779.72'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.63'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.45'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.73'

# Volume(1h) column (E) updates
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.61%  '
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("E24").Value = '  +6.82%  '
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("E26").Value = '  +2.09%  '
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  +0.45%  '
